$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: cardholder name / card number / surname ---
$ws.Range("C2").Value = "Hartmut"
# Card number must stay text (it's a 16-digit string, not a numeric value),
# so force the text number format before writing it - otherwise Excel's
# auto-detection would turn the digit string into a Number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 22.07.2024"

# --- Row 6 (was the Allianz entry, now the electricity entry) ---
$ws.Range("B6").Value = "24.07."
$ws.Range("C6").Value = "25.07."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 1506622"
$ws.Range("E6").Value = "86,14-"

# --- Row 7 (was the electricity entry, now the Allianz entry) ---
$ws.Range("B7").Value = "27.07."
$ws.Range("C7").Value = "28.07."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-14761576"
$ws.Range("E7").Value = "57,15-"

# --- Row 8 (was PayPal, now Amazon) ---
$ws.Range("B8").Value = "28.07."
$ws.Range("C8").Value = "29.07."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU PDSAMG"
$ws.Range("E8").Value = "49,54-"

# --- Rows 9, 10 and 11 are emptied out (no more transactions listed) ---
# Note: the E column cells are merged (E:F), so plain .ClearContents() on
# those can silently no-op against the merged range in this engine -
# assigning an empty string is the reliable way to blank them out.
$ws.Range("B9:D9").Value = ""
$ws.Range("B10:D10").Value = ""
$ws.Range("B11:D11").Value = ""

$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

$ws.Range("E11").Value = ""
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 01.08.2024"
$ws.Range("E12").Value = "192,83-"

# --- Next statement date footer ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 09.08.2024"
